# Add "Swap Nodes in Pairs" (LeetCode 24) as a new row (58) to the
# LeetCode problem tracker worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the new row's plain data columns -----------------------------
$ws.Range("A58").Value = "Swap Nodes in Pairs"
$ws.Range("B58").Value = "Linked List"
$ws.Range("C58").Value = "No"
$ws.Range("D58").Value = "Yes"
$ws.Range("E58").Value = "Medium"
$ws.Range("F58").Value = "Medium"

# --- Add the hyperlink cell in column G ------------------------------------
# Hyperlinks.Add stamps its own (duplicate) hyperlink font style onto the
# cell, so re-apply the existing Hyperlink-styled formatting from the row
# above (G57) afterwards to keep the same cell style index.
$ws.Hyperlinks.Add($ws.Range("G58"), "24 - Swap Nodes in Pairs", "", "", "24 - Swap Nodes in Pairs")
$ws.Range("G57").Copy()
$ws.Range("G58").PasteSpecial(-4122)

# --- Match the trailing selection / active cell left by the author --------
[void]$ws.Range("M62").Select()
